$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1835205992509363
$ws.Range("C2").Value = 0.5805243445692884
$ws.Range("J2").Value = 0.01872659176029963
$ws.Range("P2").Value = 0.1310861423220974
$ws.Range("S2").Value = 0.08614232209737828
$ws.Range("C3").Value = 0.01875
$ws.Range("J3").Value = 0.0625
$ws.Range("P3").Value = 0.6875
$ws.Range("S3").Value = 0.23125
$ws.Range("P4").Value = 0.8095238095238095
$ws.Range("S4").Value = 0.1904761904761905
$ws.Range("B6").Value = 0.07936507936507936
$ws.Range("F6").Value = 0.05291005291005291
$ws.Range("J6").Value = 0.2804232804232804
$ws.Range("O6").Value = 0.02116402116402116
$ws.Range("Q6").Value = 0.1058201058201058
$ws.Range("R6").Value = 0.08994708994708994
$ws.Range("S6").Value = 0.3703703703703703
$ws.Range("B7").Value = 0.07589285714285714
$ws.Range("D7").Value = 0.01785714285714286
$ws.Range("F7").Value = 0.03571428571428571
$ws.Range("J7").Value = 0.1651785714285714
$ws.Range("O7").Value = 0.01785714285714286
$ws.Range("Q7").Value = 0.1607142857142857
$ws.Range("R7").Value = 0.08482142857142858
$ws.Range("S7").Value = 0.4419642857142857
$ws.Range("B8").Value = 0.08181818181818182
$ws.Range("D8").Value = 0.00909090909090909
$ws.Range("F8").Value = 0.03636363636363636
$ws.Range("J8").Value = 0.1151515151515152
$ws.Range("O8").Value = 0.00909090909090909
$ws.Range("Q8").Value = 0.1454545454545454
$ws.Range("R8").Value = 0.06666666666666667
$ws.Range("S8").Value = 0.5363636363636364
$ws.Range("B9").Value = 0.08205128205128205
$ws.Range("D9").Value = 0.02564102564102564
$ws.Range("F9").Value = 0.04102564102564103
$ws.Range("J9").Value = 0.1743589743589744
$ws.Range("O9").Value = 0.01538461538461539
$ws.Range("Q9").Value = 0.1897435897435897
$ws.Range("R9").Value = 0.1025641025641026
$ws.Range("S9").Value = 0.3692307692307693
$ws.Range("B10").Value = 0.09781021897810219
$ws.Range("D10").Value = 0.02408759124087591
$ws.Range("E10").Value = 0.00145985401459854
$ws.Range("F10").Value = 0.06788321167883211
$ws.Range("J10").Value = 0.145985401459854
$ws.Range("O10").Value = 0.01605839416058394
$ws.Range("Q10").Value = 0.1985401459854015
$ws.Range("R10").Value = 0.06642335766423357
$ws.Range("S10").Value = 0.3817518248175182
$ws.Range("G11").Value = 0.1400560224089636
$ws.Range("J11").Value = 0.1120448179271709
$ws.Range("K11").Value = 0.2016806722689076
$ws.Range("L11").Value = 0.5406162464985994
$ws.Range("S11").Value = 0.005602240896358543
$ws.Range("G12").Value = 0.774869109947644
$ws.Range("J12").Value = 0.193717277486911
$ws.Range("L12").Value = 0.005235602094240838
$ws.Range("S12").Value = 0.02617801047120419
$ws.Range("G13").Value = 0.6829268292682927
$ws.Range("J13").Value = 0.2926829268292683
$ws.Range("S13").Value = 0.02439024390243903
$ws.Range("G14").Value = 0.75
$ws.Range("J14").Value = 0.25
$ws.Range("F15").Value = 0.005076142131979695
$ws.Range("H15").Value = 0.1319796954314721
$ws.Range("I15").Value = 0.116751269035533
$ws.Range("J15").Value = 0.3604060913705584
$ws.Range("K15").Value = 0.07106598984771574
$ws.Range("M15").Value = 0.01015228426395939
$ws.Range("N15").Value = 0.005076142131979695
$ws.Range("O15").Value = 0.05076142131979695
$ws.Range("S15").Value = 0.2487309644670051
$ws.Range("F16").Value = 0.01704545454545454
$ws.Range("H16").Value = 0.1420454545454546
$ws.Range("I16").Value = 0.07954545454545454
$ws.Range("J16").Value = 0.4772727272727273
$ws.Range("K16").Value = 0.1193181818181818
$ws.Range("M16").Value = 0.005681818181818182
$ws.Range("O16").Value = 0.02272727272727273
$ws.Range("S16").Value = 0.1363636363636364
$ws.Range("F17").Value = 0.02205882352941177
$ws.Range("H17").Value = 0.09068627450980392
$ws.Range("I17").Value = 0.09558823529411764
$ws.Range("J17").Value = 0.482843137254902
$ws.Range("K17").Value = 0.1372549019607843
$ws.Range("M17").Value = 0.01470588235294118
$ws.Range("N17").Value = 0.007352941176470588
$ws.Range("O17").Value = 0.04656862745098039
$ws.Range("S17").Value = 0.1029411764705882
$ws.Range("F18").Value = 0.01796407185628742
$ws.Range("H18").Value = 0.1197604790419162
$ws.Range("I18").Value = 0.1017964071856287
$ws.Range("J18").Value = 0.4730538922155689
$ws.Range("K18").Value = 0.125748502994012
$ws.Range("M18").Value = 0.01197604790419162
$ws.Range("O18").Value = 0.04790419161676647
$ws.Range("S18").Value = 0.1017964071856287
$ws.Range("F19").Value = 0.01238390092879257
$ws.Range("H19").Value = 0.1718266253869969
$ws.Range("I19").Value = 0.07972136222910217
$ws.Range("J19").Value = 0.3893188854489164
$ws.Range("K19").Value = 0.130030959752322
$ws.Range("M19").Value = 0.02476780185758514
$ws.Range("O19").Value = 0.07662538699690402
$ws.Range("S19").Value = 0.1153250773993808
